$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update display_id values to include the {key} placeholder before processing
$ws.Range("A2").Value = "c1_{key}"
$ws.Range("A3").Value = "c2_{key}"
$ws.Range("A4").Value = "c3_{key}"
